$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.143.41"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.791.77"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'228.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").Value = "'0.556"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.78%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'31.56"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").Value = "'46.27"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("E11").Value = "  -2.20%  "
$ws.Range("D12").Value = "'0.0929"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "2.048.93"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "'11.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +12.12%  "
$ws.Range("D15").Value = "1.785.44"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").Value = "'0.640"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "34.135.81"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "'4.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("D19").Value = "'69.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "'253.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'10.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("D26").Value = "'157.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.28%  "
$ws.Range("D27").Value = "'16.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").Value = "'7.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("D29").Value = "'0.114"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -1.27%  "
$ws.Range("E34").Value = "  +1.05%  "
$ws.Range("E35").Value = "  +1.53%  "
$ws.Range("D36").Value = "1.457.86"
$ws.Range("E36").Value = "  -8.20%  "
$ws.Range("D37").Value = "'1.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").Value = "'0.632"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "'0.0188"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").Value = "'83.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").Value = "'0.909"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("D45").Value = "'0.0512"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.91%  "
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("D47").Value = "1.949.11"
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("D48").Value = "'5.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'12.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.85%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "'1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "'51.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.62%  "
